# Recording-utility update: add a new user row (Test Manager) to the Users sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

# New user record in row 4
$ws.Range("A4").Value = "test@click2cloud.com"
$ws.Range("B4").Value = "'123"
$ws.Range("C4").Value = "Test Manager"

# Extend the dropdown (list) validation on column C down to the new row
$ws.Range("C2:C4").Validation.Add(3, 1, 1, "=Sheet3!`$C`$4:`$C`$10") | Out-Null

# Resize column C (drop the old "best fit" auto width) and update the active cell
$ws.Columns.Item(3).ColumnWidth = 12.835
$ws.Range("B13").Select() | Out-Null
